# Refresh cached Universalis market-data snapshot (currentAveragePrice*,
# LevePrice*/LeveProfit*) for each crafting job's Leve-profit table, as
# produced by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 821.2727
$ws.Range("I38").Value = 470
$ws.Range("J38").Value = 899.3333
$ws.Range("K38").Value = 1410
$ws.Range("L38").Value = 2697.9999
$ws.Range("M38").Value = -1038
$ws.Range("N38").Value = -3441.9999
$ws.Range("H39").Value = 120.75
$ws.Range("I39").Value = 64.42856999999999
$ws.Range("J39").Value = 199.6
$ws.Range("K39").Value = 193.28571
$ws.Range("L39").Value = 598.8
$ws.Range("M39").Value = 102.71429
$ws.Range("N39").Value = -1190.8
$ws.Range("H88").Value = 1700
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1700
$ws.Range("K88").Value = 0
$ws.Range("N88").Value = -2512
$ws.Range("H91").Value = 1700
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1700
$ws.Range("K91").Value = 0
$ws.Range("N91").Value = -4508
$ws.Range("H92").Value = 804.7059
$ws.Range("I92").Value = 795.38464
$ws.Range("J92").Value = 835
$ws.Range("K92").Value = 795.38464
$ws.Range("L92").Value = 835
$ws.Range("M92").Value = 452.61536
$ws.Range("N92").Value = -3331
$ws.Range("H96").Value = 20834298
$ws.Range("I96").Value = 27778566
$ws.Range("J96").Value = 1493.3334
$ws.Range("K96").Value = 83335698
$ws.Range("L96").Value = 4480.0002
$ws.Range("M96").Value = -83334325
$ws.Range("N96").Value = -7226.0002
$ws.Range("H100").Value = 2635.4119
$ws.Range("I100").Value = 2233.6667
$ws.Range("J100").Value = 3599.6
$ws.Range("K100").Value = 2233.6667
$ws.Range("L100").Value = 3599.6
$ws.Range("M100").Value = -1692.6667
$ws.Range("N100").Value = -4681.6
$ws.Range("H112").Value = 3473255.2
$ws.Range("I112").Value = 566.3333
$ws.Range("J112").Value = 3832499
$ws.Range("K112").Value = 1698.9999
$ws.Range("L112").Value = 11497497
$ws.Range("M112").Value = -590.9999
$ws.Range("N112").Value = -11499713
$ws.Range("M88").ClearContents()
$ws.Range("M91").ClearContents()
$ws.Range("H125:N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 42856
$ws.Range("J7").Value = 42856
$ws.Range("L7").Value = 42856
$ws.Range("N7").Value = -43084
$ws.Range("H32").Value = 22046.98
$ws.Range("I32").Value = 22189.52
$ws.Range("K32").Value = 22189.52
$ws.Range("M32").Value = -21902.52
$ws.Range("H97").Value = 745.2381
$ws.Range("I97").Value = 803.5294
$ws.Range("J97").Value = 497.5
$ws.Range("K97").Value = 803.5294
$ws.Range("L97").Value = 497.5
$ws.Range("M97").Value = -307.5294
$ws.Range("N97").Value = -1489.5
$ws.Range("H102").Value = 1302.0714
$ws.Range("I102").Value = 1111.7273
$ws.Range("K102").Value = 1111.7273
$ws.Range("M102").Value = 510.2727
$ws.Range("H132").Value = 11603.49
$ws.Range("I132").Value = 1525.421
$ws.Range("J132").Value = 41062.46
$ws.Range("K132").Value = 4576.263
$ws.Range("L132").Value = 123187.38
$ws.Range("M132").Value = -2046.263
$ws.Range("N132").Value = -128247.38

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2357.7083
$ws.Range("I132").Value = 1179.4667
$ws.Range("K132").Value = 3538.4001
$ws.Range("M132").Value = -1008.4001
$ws.Range("H134").Value = 1098.3793
$ws.Range("I134").Value = 906.0625
$ws.Range("K134").Value = 2718.1875
$ws.Range("M134").Value = -183.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4309.5454
$ws.Range("J81").Value = 4440.5
$ws.Range("L81").Value = 13321.5
$ws.Range("N81").Value = -15567.5
$ws.Range("H84").Value = 4309.5454
$ws.Range("J84").Value = 4440.5
$ws.Range("L84").Value = 39964.5
$ws.Range("N84").Value = -51196.5
$ws.Range("H131").Value = 765.22
$ws.Range("J131").Value = 782.72046
$ws.Range("L131").Value = 2348.16138
$ws.Range("N131").Value = -12428.16138
$ws.Range("H132").Value = 1272.5555
$ws.Range("I132").Value = 950.4286
$ws.Range("J132").Value = 2400
$ws.Range("K132").Value = 8553.857399999999
$ws.Range("L132").Value = 21600
$ws.Range("M132").Value = -6023.857399999999
$ws.Range("N132").Value = -26660
$ws.Range("H139").Value = 1983.6471
$ws.Range("I139").Value = 1195.1818
$ws.Range("K139").Value = 3585.5454
$ws.Range("M139").Value = 1554.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3934.0833
$ws.Range("J80").Value = 3988.625
$ws.Range("L80").Value = 3988.625
$ws.Range("N80").Value = -5984.625
$ws.Range("H83").Value = 3934.0833
$ws.Range("J83").Value = 3988.625
$ws.Range("L83").Value = 19943.125
$ws.Range("N83").Value = -29927.125
$ws.Range("H132").Value = 62662.76
$ws.Range("I132").Value = 55182.684
$ws.Range("J132").Value = 86349.664
$ws.Range("K132").Value = 165548.052
$ws.Range("L132").Value = 259048.992
$ws.Range("M132").Value = -163018.052
$ws.Range("N132").Value = -264108.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 313.21054
$ws.Range("I16").Value = 312.85715
$ws.Range("J16").Value = 314.2
$ws.Range("K16").Value = 312.85715
$ws.Range("L16").Value = 314.2
$ws.Range("M16").Value = -142.85715
$ws.Range("N16").Value = -654.2
$ws.Range("H136").Value = 1525.8182
$ws.Range("I136").Value = 1380.4
$ws.Range("K136").Value = 4141.200000000001
$ws.Range("M136").Value = -1591.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1298.2941
$ws.Range("I132").Value = 761.1818
$ws.Range("J132").Value = 2283
$ws.Range("K132").Value = 2283.5454
$ws.Range("L132").Value = 6849
$ws.Range("M132").Value = 246.4546
$ws.Range("N132").Value = -11909
$ws.Range("H136").Value = 20834822
$ws.Range("I136").Value = 33334666
$ws.Range("K136").Value = 100003998
$ws.Range("M136").Value = -100001448
